# "office soft assertion implementation"
#
# 1) Keywords sheet: add a new "ProceedOnFail" column (F) - Y for the
#    openBrowser/navigate/validateTitle/click steps, N for the rest.
# 2) Data sheet: the LoginTest "Browser" data row for U1/P1 switches from
#    Mozilla to Chrome, and a new Edge/U3/P3/LoginSuccess row is appended
#    to the LoginTest block (duplicating the existing last row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Keywords sheet - new "ProceedOnFail" column
# ---------------------------------------------------------------------
$wsKeywords = $wb.Worksheets.Item("Keywords")
$wsKeywords.Activate()

# Header, formatted like the rest of row 1
$wsKeywords.Range("E1").Copy()
$wsKeywords.Range("F1").PasteSpecial(-4122)
$wsKeywords.Range("F1").Value = "ProceedOnFail"

# Body cells, formatted like the rest of the data rows
$wsKeywords.Range("E2").Copy()
$wsKeywords.Range("E10:F13").PasteSpecial(-4122)
$wsKeywords.Range("F2:F9").PasteSpecial(-4122)

$wsKeywords.Range("F2:F5").Value = "Y"
$wsKeywords.Range("F6:F13").Value = "N"

$wsKeywords.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

$wsKeywords.Range("E6").Select()

# ---------------------------------------------------------------------
# Data sheet - Browser data updates
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")
$wsData.Activate()

# First login test-data row now uses Chrome instead of Mozilla
$wsData.Range("B3").Value = "Chrome"

# Duplicate the Edge row as a new row 6
$wsData.Rows.Item(6).Insert()
$wsData.Range("A5:E5").Copy()
$wsData.Range("A6:E6").PasteSpecial(-4104)

$wsData.Range("A6:E6").Select()
